$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.233.40"
$ws.Range("E2").Value = "  +1.39%  "

$ws.Range("D3").Value = "1.906.79"
$ws.Range("E3").Value = "  +2.01%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'307.61"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  +3.50%  "

$ws.Range("D8").Value = "'0.3787"
$ws.Range("E8").Value = "  +3.58%  "

$ws.Range("D9").Value = "'0.07295"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("E10").Value = "  +3.05%  "

$ws.Range("D11").Value = "'0.9003"
$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07671"
$ws.Range("E12").Value = "  +1.96%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.914.24"
$ws.Range("E13").Value = "  +2.41%  "

$ws.Range("D14").Value = "'95.04"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "'5.260"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "'0.000008690"
$ws.Range("E17").Value = "  +2.29%  "

$ws.Range("D18").Value = "'14.51"

$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").Value = "27.299.68"
$ws.Range("E20").Value = "  +1.47%  "

$ws.Range("D21").Value = "'5.083"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").Value = "2.148.50"
$ws.Range("E22").Value = "  +1.67%  "

$ws.Range("D24").Value = "'6.442"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("D25").Value = "'2.312"
$ws.Range("E25").Value = "  +10.04%  "

$ws.Range("D26").Value = "'145.89"
$ws.Range("E26").Value = "  -1.45%  "

$ws.Range("D27").Value = "'18.15"
$ws.Range("E27").Value = "  +1.66%  "

$ws.Range("D28").Value = "'1.729"
$ws.Range("E28").Value = "  -3.09%  "

$ws.Range("D29").Value = "'114.92"
$ws.Range("E29").Value = "  +1.37%  "

$ws.Range("D30").Value = "'4.960"
$ws.Range("E30").Value = "  +4.93%  "

$ws.Range("D31").Value = "'4.826"
$ws.Range("E31").Value = "  +2.64%  "

$ws.Range("E32").Value = "  +0.89%  "

$ws.Range("D33").Value = "'0.05081"
$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.248"
$ws.Range("E34").Value = "  +8.10%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7946"
$ws.Range("E35").Value = "  +6.38%  "

$ws.Range("D36").Value = "'2.998"
$ws.Range("E36").Value = "  +0.51%  "

$ws.Range("D37").Value = "'3.304"

$ws.Range("D38").Value = "'2.622"
$ws.Range("E38").Value = "  +3.65%  "

$ws.Range("D39").Value = "'0.5682"
$ws.Range("E39").Value = "  +1.63%  "

$ws.Range("D40").Value = "'0.01998"
$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("D41").Value = "'1.074"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").Value = "'9.016"
$ws.Range("E42").Value = "  +4.95%  "

$ws.Range("D43").Value = "'6.652"

$ws.Range("D44").Value = "'119.19"
$ws.Range("E44").Value = "  +3.20%  "

$ws.Range("D45").Value = "'0.1522"
$ws.Range("E45").Value = "  +3.20%  "

$ws.Range("D46").Value = "'0.4871"
$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = "  +0.97%  "

$ws.Range("D48").Value = "'0.9987"
$ws.Range("E48").Value = "  -0.13%  "

$ws.Range("D49").Value = "'1.608"
$ws.Range("E49").Value = "  +2.78%  "

$ws.Range("D50").Value = "'37.50"
$ws.Range("E50").Value = "  +1.62%  "

$ws.Range("D51").Value = "'64.26"
$ws.Range("E51").Value = "  +1.99%  "
